$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the seller name (shared string "GUSTAVO GLUNZ" -> "GUSTAVO")
$ws.Range("B2").Value = "GUSTAVO"

# Update phone numbers in A2 and D2
$ws.Range("A2").Value = 5491161405589
$ws.Range("D2").Value = 5491161405589

# Update the active selection on the sheet to B3
$ws.Range("B3").Select()
